$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the old "Saldo" header that lived in E1 on the main sheet;
# a dedicated "Saldo" sheet is being introduced instead (see below).
$ws.Range("E1").Value = $null

# Add two sample transaction rows to the main "Sheet"
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "123/321"
$ws.Range("C2").Value = "'05/05/2004"
$ws.Range("C2").ClearFormats()
$ws.Range("D2").Value = "compras"

$ws.Range("A3").Value = 5
$ws.Range("B3").Value = "123/321"
$ws.Range("C3").Value = "'05/05/2004"
$ws.Range("C3").ClearFormats()
$ws.Range("D3").Value = "compras"

# Add a new "Saldo" worksheet after the existing one, to start
# implementing the Transferir functionality.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$saldo = $wb.Worksheets.Add($null, $lastSheet)
$saldo.Name = "Saldo"

$saldo.Range("A1").Value = "Saldo"
$ws.Range("A1").Copy()
$saldo.Range("A1").PasteSpecial(-4122)

$saldo.Range("A2").Value = 0
$saldo.Range("A3").Value = 10
$saldo.Range("A4").Value = 5

$saldo.Range("A3").Select() | Out-Null
